$d = $word.ActiveDocument

$replacements = @(
    @("75÷8=9, 3", "77÷6=12, 5"),
    @("53÷8=6, 5", "35÷8=4, 3"),
    @("81÷4=20, 1", "26÷5=5, 1"),
    @("82÷5=16, 2", "83÷6=13, 5"),
    @("27÷6=4, 3", "40÷8=5, 0"),
    @("95÷2=47, 1", "39÷3=13, 0"),
    @("77÷5=15, 2", "85÷5=17, 0"),
    @("48÷9=5, 3", "76÷7=10, 6"),
    @("42÷8=5, 2", "13÷3=4, 1"),
    @("45÷2=22, 1", "17÷4=4, 1"),
    @("24÷3=8, 0", "15÷9=1, 6"),
    @("27÷7=3, 6", "51÷8=6, 3"),
    @("39÷2=19, 1", "94÷9=10, 4"),
    @("82÷7=11, 5", "61÷9=6, 7"),
    @("42÷2=21, 0", "70÷7=10, 0"),
    @("69÷4=17, 1", "18÷2=9, 0"),
    @("67÷3=22, 1", "90÷7=12, 6"),
    @("65÷2=32, 1", "75÷3=25, 0"),
    @("41÷9=4, 5", "92÷7=13, 1"),
    @("51÷7=7, 2", "71÷2=35, 1"),
    @("23÷6=3, 5", "82÷7=11, 5"),
    @("14÷7=2, 0", "55÷9=6, 1"),
    @("43÷3=14, 1", "28÷9=3, 1"),
    @("82÷9=9, 1", "62÷5=12, 2"),
    @("77÷7=11, 0", "62÷8=7, 6")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
